{"js": "// [oldText, newText] pairs, in document order. Every oldText value is\n// unique within this document, so a single (first-hit) search+replace per\n// pair is unambiguous regardless of processing order.\nconst replacements = [\n  [\"2025-07-22 Tuesday\", \"2025-07-23 Wednesday\"],\n  [\"45\u00d793=4185\", \"87\u00d752=4524\"],\n  [\"39\u00d779=3081\", \"85\u00d787=7395\"],\n  [\"34\u00d780=2720\", \"20\u00d794=1880\"],\n  [\"27\u00d785=2295\", \"98\u00d779=7742\"],\n  [\"96\u00d766=6336\", \"30\u00d728=840\"],\n  [\"45\u00d723=1035\", \"83\u00d714=1162\"],\n  [\"79\u00d789=7031\", \"73\u00d756=4088\"],\n  [\"60\u00d798=5880\", \"63\u00d714=882\"],\n  [\"79\u00d790=7110\", \"18\u00d794=1692\"],\n  [\"61\u00d745=2745\", \"89\u00d745=4005\"],\n  [\"55\u00d735=1925\", \"13\u00d717=221\"],\n  [\"51\u00d768=3468\", \"36\u00d770=2520\"],\n  [\"88\u00d724=2112\", \"12\u00d737=444\"],\n  [\"79\u00d716=1264\", \"87\u00d752=4524\"],\n  [\"33\u00d718=594\", \"31\u00d722=682\"],\n  [\"66\u00d763=4158\", \"36\u00d754=1944\"],\n  [\"12\u00d770=840\", \"82\u00d780=6560\"],\n  [\"29\u00d787=2523\", \"73\u00d781=5913\"],\n  [\"77\u00d756=4312\", \"80\u00d719=1520\"],\n  [\"25\u00d715=375\", \"55\u00d798=5390\"],\n  [\"17\u00d760=1020\", \"50\u00d759=2950\"],\n  [\"41\u00d784=3444\", \"36\u00d784=3024\"],\n  [\"58\u00d771=4118\", \"50\u00d711=550\"],\n  [\"72\u00d732=2304\", \"93\u00d716=1488\"],\n  [\"78\u00d784=6552\", \"77\u00d737=2849\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  // Replace only the first occurrence (sources are unique in this document).\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# [oldText, newText] pairs, in document order. Every oldText value is unique\n# within this document, so a single wdReplaceOne pass per pair is\n# unambiguous regardless of processing order.\n$replacements = @(\n    @(\"2025-07-22 Tuesday\", \"2025-07-23 Wednesday\"),\n    @(\"45\u00d793=4185\", \"87\u00d752=4524\"),\n    @(\"39\u00d779=3081\", \"85\u00d787=7395\"),\n    @(\"34\u00d780=2720\", \"20\u00d794=1880\"),\n    @(\"27\u00d785=2295\", \"98\u00d779=7742\"),\n    @(\"96\u00d766=6336\", \"30\u00d728=840\"),\n    @(\"45\u00d723=1035\", \"83\u00d714=1162\"),\n    @(\"79\u00d789=7031\", \"73\u00d756=4088\"),\n    @(\"60\u00d798=5880\", \"63\u00d714=882\"),\n    @(\"79\u00d790=7110\", \"18\u00d794=1692\"),\n    @(\"61\u00d745=2745\", \"89\u00d745=4005\"),\n    @(\"55\u00d735=1925\", \"13\u00d717=221\"),\n    @(\"51\u00d768=3468\", \"36\u00d770=2520\"),\n    @(\"88\u00d724=2112\", \"12\u00d737=444\"),\n    @(\"79\u00d716=1264\", \"87\u00d752=4524\"),\n    @(\"33\u00d718=594\", \"31\u00d722=682\"),\n    @(\"66\u00d763=4158\", \"36\u00d754=1944\"),\n    @(\"12\u00d770=840\", \"82\u00d780=6560\"),\n    @(\"29\u00d787=2523\", \"73\u00d781=5913\"),\n    @(\"77\u00d756=4312\", \"80\u00d719=1520\"),\n    @(\"25\u00d715=375\", \"55\u00d798=5390\"),\n    @(\"17\u00d760=1020\", \"50\u00d759=2950\"),\n    @(\"41\u00d784=3444\", \"36\u00d784=3024\"),\n    @(\"58\u00d771=4118\", \"50\u00d711=550\"),\n    @(\"72\u00d732=2304\", \"93\u00d716=1488\"),\n    @(\"78\u00d784=6552\", \"77\u00d737=2849\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # wdFindStop = 0 (no wrap), wdReplaceOne = 1 (replace the single match)\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 0, $false, $newText, 1)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
